# Working Apollo - description removed
#
# 1) Append a "." (its own run) right after "... final push".
# 2) Fill in the previously-empty paragraph that follows with the new
#    21.2 diary entry.
# 3) Remove the "-v praxi otestovana ..." bullet under LESSONS LEARNED
#    (whole paragraph, including its paragraph mark).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: find a (top-level body) paragraph whose text starts with
# $prefix - lets us locate anchors without hard-coding paragraph
# indices.
# ---------------------------------------------------------------------
function Get-ParagraphStartingWith($doc, $prefix) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.StartsWith($prefix)) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Locate the "15.2 ... final push" paragraph and, inside it, the
# already-formatted "push" run (sz/szCs = 26) - it is used below as a
# formatting "seed" so newly-inserted text keeps matching run
# properties instead of falling back to un-sized default formatting.
# ---------------------------------------------------------------------
$pushPara = Get-ParagraphStartingWith $d "15.2"
$pushRange = $pushPara.Range
$seedSource = $d.Range($pushRange.End - 5, $pushRange.End - 1)   # the word "push"

# ---------------------------------------------------------------------
# Step 1: add a new "." run at the end of the "final push" paragraph.
# Seeding + overwrite (rather than plain InsertAfter) means the new
# run picks up the full rPr (sz AND szCs), not just sz.
# ---------------------------------------------------------------------
$insertPoint = $d.Range($pushRange.End - 1, $pushRange.End - 1)
$insertPoint.FormattedText = $seedSource.FormattedText

$pushPara = Get-ParagraphStartingWith $d "15.2"
$pushRange = $pushPara.Range
$newRunRange = $d.Range($pushRange.End - 5, $pushRange.End - 1)
$newRunRange.Text = "."

# ---------------------------------------------------------------------
# Step 2: fill the empty paragraph that follows with the 21.2 entry,
# using the same seed/overwrite trick (that paragraph starts out with
# zero runs, so there is nothing to Find/Replace into).
# ---------------------------------------------------------------------
$pushPara = Get-ParagraphStartingWith $d "15.2"
$emptyPara = $pushPara.Next()
$emptyRange = $emptyPara.Range
$insertPoint2 = $d.Range($emptyRange.End - 1, $emptyRange.End - 1)
$insertPoint2.FormattedText = $seedSource.FormattedText

$pushPara = Get-ParagraphStartingWith $d "15.2"
$emptyPara = $pushPara.Next()
$emptyRange = $emptyPara.Range
$contentRange = $d.Range($emptyRange.Start, $emptyRange.End - 1)
$contentRange.Text = "21.2 – Opravené Apollo po konzultaci se spolužákem. Stačilo odstranit description v GraphTypeDefinitions u GroupGQLModel, který je tahán z jiného projektu. Až teď opravdu dokončený projekt."

# ---------------------------------------------------------------------
# Step 3: drop the "-v praxi otestovana ..." bullet entirely (text +
# its paragraph mark), leaving LESSONS LEARNED followed directly by
# "-první cesta ...".
# ---------------------------------------------------------------------
$bulletPara = Get-ParagraphStartingWith $d "-v praxi"
$bulletPara.Range.Delete() | Out-Null
